$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values for each row, re-pulled from source data.
$values = @{
    2  = 3
    3  = -1
    4  = -3
    5  = 3
    6  = 2
    7  = 6
    8  = -2
    9  = 3
    10 = -6
    11 = -2
    12 = -4
    13 = 7
    14 = 1
    15 = 1
    16 = -1
    17 = -4
    18 = -1
    20 = 7
    21 = 3
    22 = 6
    23 = 1
    24 = 3
    25 = -1
    26 = -6
    29 = 1
    30 = 7
    31 = 1
    33 = -4
    34 = -5
    35 = 8
    37 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
